# Insert a new data row at row 19 (pushing the existing rows 19-98 down to
# 20-99, which also pushes the former last row, old row 98, down to the new
# last row 99). Then populate the newly inserted row 19 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 5
$ws.Cells.Item(19, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(19, 3).Value = 'Maule'
$ws.Cells.Item(19, 4).Value = 44623
$ws.Cells.Item(19, 5).Value = 7
$ws.Cells.Item(19, 6).Value = 'Fruta'
$ws.Cells.Item(19, 7).Value = 100103
$ws.Cells.Item(19, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(19, 9).Value = 100103002
$ws.Cells.Item(19, 10).Value = 'Ciruela'
$ws.Cells.Item(19, 11).Value = 'Black Amber'
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 500
$ws.Cells.Item(19, 14).Value = 9000
$ws.Cells.Item(19, 15).Value = 9000
$ws.Cells.Item(19, 16).Value = 9000
$ws.Cells.Item(19, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 500
$ws.Cells.Item(19, 20).Value = 18
